$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns (AD, AE, AF) on row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting used by the other header cells in row 1
# (bold font, thin box border, centered / top-aligned) by copying the
# existing header cell format onto the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill team record values (Wins=89, Losses=73, Ties=0) for every player row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 89  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 73  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
